$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.75
$ws.Range("I2").Value = 1.75
$ws.Range("K2").Value = 2.3
$ws.Range("L2").Value = 2.38
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("AB2").Value = 34
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 15
$ws.Range("AH2").Value = 8.5
$ws.Range("AI2").Value = 9.5
$ws.Range("AT2").Value = 3.25
